$wb = $excel.ActiveWorkbook

# 1. Replace the shared string "Ready for handoff" with "In Translation"
#    wherever it appears across all worksheets (Overview zh-cn/de-de status
#    columns, and the Status column of each per-language sheet).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value()
        if ("Ready for handoff" -eq $val) {
            $cell.Value = "In Translation"
        }
    }
}

# 2. Narrow the "Status" related columns (E & F on Overview, C on the
#    per-language sheets) to match the regenerated report layout.
#    ColumnWidth is quantized internally to 1/6-character increments, so
#    12.5 is the input that resolves to the closest achievable width.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").ColumnWidth = 12.5
